$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Add the new "Gardener" ally unit strings below the existing "Mallet Guy" rows
$ws.Range("A20").Value = "unitAllyGardener"
$ws.Range("B20").Value = "Gardener"
$ws.Range("A21").Value = "unitAllyGardenerDesc"
$ws.Range("B21").Value = "He gardens stuff."

# Move the active selection to reflect the next empty row, like the diff shows
$ws.Range("A22").Select()
